# Update cryptos list values (price + 1h volume change) per latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.523.08"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "1.672.17"
$ws.Range("E3").Value = "  +2.35%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "219.61"
$ws.Range("E5").Value = "  +2.45%  "
$ws.Range("D6").Value = "'0.530"
$ws.Range("E6").Value = "  +2.35%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "29.56"
$ws.Range("E8").Value = "  +3.53%  "
$ws.Range("E9").Value = "  +2.77%  "
$ws.Range("D10").Value = "0.0642"
$ws.Range("E10").Value = "  +5.63%  "
$ws.Range("D11").Value = "0.0905"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "1.914.06"
$ws.Range("E12").Value = "  +2.47%  "
$ws.Range("D13").Value = "0.614"
$ws.Range("E13").Value = "  +8.91%  "
$ws.Range("D14").Value = "1.668.99"
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("D15").Value = "10.18"
$ws.Range("E15").Value = "  +9.08%  "
$ws.Range("E16").Value = "  +3.72%  "
$ws.Range("D17").Value = "30.547.38"
$ws.Range("E17").Value = "  +2.03%  "
$ws.Range("D18").Value = "66.48"
$ws.Range("E18").Value = "  +3.76%  "
$ws.Range("D19").Value = "242.66"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Value = "0.0₃0722"
$ws.Range("E20").Value = "  +3.12%  "
$ws.Range("D21").Value = "0.998"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("E22").Value = "  +3.31%  "
$ws.Range("D23").Value = "9.98"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").Value = "158.64"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("D26").Value = "15.83"
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("D27").Value = "0.112"
$ws.Range("E27").Value = "  +2.33%  "
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "0.0496"
$ws.Range("E30").Value = "  +1.83%  "
$ws.Range("E31").Value = "  +2.81%  "
$ws.Range("E32").Value = "  +2.63%  "
$ws.Range("E33").Value = "  +3.09%  "
$ws.Range("D34").Value = "1.497.04"
$ws.Range("E34").Value = "  +4.99%  "
$ws.Range("E35").Value = "  +7.12%  "
$ws.Range("D36").Value = "84.18"
$ws.Range("E36").Value = "  +10.62%  "
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("D38").Value = "0.601"
$ws.Range("E38").Value = "  +8.72%  "
$ws.Range("E39").Value = "  +4.99%  "
$ws.Range("E40").Value = "  -4.47%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("E42").Value = "  +1.10%  "
$ws.Range("E43").Value = "  +1.74%  "
$ws.Range("E44").Value = "  -1.84%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Value = "5.54"
$ws.Range("E47").Value = "  +3.21%  "
$ws.Range("D48").Value = "50.84"
$ws.Range("E48").Value = "  -4.00%  "
$ws.Range("D49").Value = "1.806.28"
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("D50").Value = "94.87"
$ws.Range("E50").Value = "  +4.37%  "
$ws.Range("E51").Value = "  -0.10%  "
